# Avances Etiquetado Roboflow 6/5/2025
# Fill in this week's progress row (row 30, date 06/05/2025) with the
# latest counts and mark the "Notas" column as N/A, same as the prior weeks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E30").Value = 127
$ws.Range("F30").Value = 234
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 1012
$ws.Range("J30").Value = "N/A"

# Scroll the view down / move the selection, mirroring the author's
# last on-screen position when saving.
$excel.ActiveWindow.ScrollRow = 18
$ws.Range("J33").Select()
